# Auto-generated edit script: apply 2022-11-12 crime data updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 6376
$ws.Range("I3").Value = 6662
$ws.Range("I4").Value = 1521
$ws.Range("I5").Value = 616
$ws.Range("I6").Value = 7572
$ws.Range("I7").Value = 22747

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 213
$ws.Range("I7").Value = 715

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I4").Value = 28
$ws.Range("I6").Value = 105
$ws.Range("I7").Value = 410

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 208
$ws.Range("I3").Value = 324
$ws.Range("I6").Value = 264
$ws.Range("I7").Value = 875

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I2").Value = 73
$ws.Range("I7").Value = 198

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 179
$ws.Range("I7").Value = 723
$ws.Range("I8").Value = 1372
$ws.Range("I10").Value = 160
$ws.Range("I11").Value = 341
$ws.Range("I12").Value = 54
$ws.Range("I13").Value = 39
$ws.Range("I15").Value = 263
$ws.Range("I19").Value = 639
$ws.Range("I20").Value = 564
$ws.Range("I21").Value = 101
$ws.Range("I23").Value = 223
$ws.Range("I26").Value = 30
$ws.Range("I29").Value = 1382
$ws.Range("I33").Value = 1029
$ws.Range("I35").Value = 31
$ws.Range("I37").Value = 715
$ws.Range("I41").Value = 97
$ws.Range("I42").Value = 811
$ws.Range("I43").Value = 197
$ws.Range("I44").Value = 169
$ws.Range("I50").Value = 113
$ws.Range("I51").Value = 272
$ws.Range("I52").Value = 488
$ws.Range("I53").Value = 247
$ws.Range("I54").Value = 464
$ws.Range("I55").Value = 257
$ws.Range("I58").Value = 13
$ws.Range("I60").Value = 126
$ws.Range("I63").Value = 71
$ws.Range("I64").Value = 186
$ws.Range("I67").Value = 875
$ws.Range("I73").Value = 209
$ws.Range("I74").Value = 36
$ws.Range("I76").Value = 327
$ws.Range("I78").Value = 310
$ws.Range("I79").Value = 648
$ws.Range("I84").Value = 198
$ws.Range("I85").Value = 1023
$ws.Range("I86").Value = 143
$ws.Range("I88").Value = 210
$ws.Range("I90").Value = 292
$ws.Range("I94").Value = 233
$ws.Range("I97").Value = 188
$ws.Range("I98").Value = 160
$ws.Range("I99").Value = 410
$ws.Range("I101").Value = 22747

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I3").Value = 378
$ws.Range("I6").Value = 330
$ws.Range("I7").Value = 1029

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I6").Value = 223
$ws.Range("I7").Value = 464

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I3").Value = 478
$ws.Range("I6").Value = 382
$ws.Range("I7").Value = 1382

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I3").Value = 189
$ws.Range("I6").Value = 200
$ws.Range("I7").Value = 639

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("I6").Value = 48
$ws.Range("I7").Value = 169

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I6").Value = 149
$ws.Range("I7").Value = 327

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 291
$ws.Range("I6").Value = 258
$ws.Range("I7").Value = 1023

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("I2").Value = 30
$ws.Range("I7").Value = 97

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 200
$ws.Range("I3").Value = 250
$ws.Range("I6").Value = 279
$ws.Range("I7").Value = 811

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range("I5").Value = 14
$ws.Range("I6").Value = 39

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("I6").Value = 73
$ws.Range("I7").Value = 160

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I6").Value = 113
$ws.Range("I7").Value = 310

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I3").Value = 82
$ws.Range("I6").Value = 80
$ws.Range("I7").Value = 257

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I3").Value = 78
$ws.Range("I7").Value = 223

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("I6").Value = 77
$ws.Range("I7").Value = 101

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I3").Value = 210
$ws.Range("I6").Value = 189
$ws.Range("I7").Value = 648

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I3").Value = 53
$ws.Range("I7").Value = 186

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I2").Value = 156
$ws.Range("I5").Value = 16
$ws.Range("I6").Value = 196
$ws.Range("I7").Value = 564

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I2").Value = 131
$ws.Range("I3").Value = 170
$ws.Range("I7").Value = 488

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I3").Value = 38
$ws.Range("I7").Value = 233

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I6").Value = 100
$ws.Range("I7").Value = 263

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I6").Value = 103
$ws.Range("I7").Value = 160

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I2").Value = 34
$ws.Range("I7").Value = 113

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 30

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I6").Value = 91
$ws.Range("I7").Value = 341

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("I6").Value = 12
$ws.Range("I7").Value = 31

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I3").Value = 66
$ws.Range("I7").Value = 209

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I6").Value = 37
$ws.Range("I7").Value = 179

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("I6").Value = 123
$ws.Range("I7").Value = 188

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I6").Value = 64
$ws.Range("I7").Value = 210

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I3").Value = 394
$ws.Range("I6").Value = 442
$ws.Range("I7").Value = 1372

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("I3").Value = 13
$ws.Range("I7").Value = 143

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I2").Value = 93
$ws.Range("I6").Value = 104
$ws.Range("I7").Value = 292

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I6").Value = 109
$ws.Range("I7").Value = 272

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I2").Value = 43
$ws.Range("I7").Value = 126

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I2").Value = 39
$ws.Range("I6").Value = 112
$ws.Range("I7").Value = 197

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I3").Value = 55
$ws.Range("I7").Value = 247

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 238
$ws.Range("I3").Value = 223
$ws.Range("I6").Value = 192
$ws.Range("I7").Value = 723

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("I6").Value = 29
$ws.Range("I7").Value = 54

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("I6").Value = 25
$ws.Range("I7").Value = 36

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("I3").Value = 3
$ws.Range("I7").Value = 13
